$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the calibration data in B7:D9, clearing the old Times New Roman
# style so the cells fall back to the default (unstyled) format.
$ws.Range("B7:D9").ClearFormats()

$ws.Range("B7").Value = 0.73478960000000004
$ws.Range("C7").Value = 1.045563
$ws.Range("D7").Value = 1.4721630000000001

$ws.Range("B8").Value = 1.0963676
$ws.Range("C8").Value = 1.0839559999999999
$ws.Range("D8").Value = 1.207802

$ws.Range("B9").Value = 1.4541297
$ws.Range("C9").Value = 1.2884869999999999
$ws.Range("D9").Value = 1.0999989999999999

# Update the active selection on the sheet.
$ws.Range("H2:H10").Select()

# Reposition the workbook window.
$excel.ActiveWindow.Left = -26720
